# Natmi following Dr Hou advice
# Replace the single-sender (sCs) LR-pair results (rows 2-4) with the
# updated two-sender (FAPs, sCs) results against all three target
# clusters (ECs, FAPs, sCs) - rows 2-7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("FAPs", "Wnt6", "Fzd7", "ECs", 2, 0.6666666666666666, 0.06347066666666668, 0.190412, 0.03928106159455794, 0.03928106159455794, 3, 1, 2.418393, 7.255179, 0.0919828589765645, 0.0919828589765645, 0.153497015972, 1.381473143748, 0.003613184349101967, 0.003613184349101966),
  @("FAPs", "Wnt6", "Fzd7", "FAPs", 2, 0.6666666666666666, 0.06347066666666668, 0.190412, 0.03928106159455794, 0.03928106159455794, 3, 1, 10.11799233333333, 30.353977, 0.3848348311969811, 0.3848348311969811, 0.642195718724889, 5.779761468524001, 0.01511672070797992, 0.01511672070797992),
  @("FAPs", "Wnt6", "Fzd7", "sCs", 2, 0.6666666666666666, 0.06347066666666668, 0.190412, 0.03928106159455794, 0.03928106159455794, 3, 1, 13.75539366666667, 41.266181, 0.5231823098264544, 0.5231823098264544, 0.873064006285778, 7.857576056572002, 0.02055115653747605, 0.02055115653747605),
  @("sCs", "Wnt6", "Fzd7", "ECs", 3, 1, 1.552337666666667, 4.657013, 0.9607189384054421, 0.960718938405442, 3, 1, 2.418393, 7.255179, 0.0919828589765645, 0.0919828589765645, 3.754162546703, 33.787462920327, 0.08836967462746254, 0.08836967462746252),
  @("sCs", "Wnt6", "Fzd7", "FAPs", 3, 1, 1.552337666666667, 4.657013, 0.9607189384054421, 0.960718938405442, 3, 1, 10.11799233333333, 30.353977, 0.3848348311969811, 0.3848348311969811, 15.70654061007789, 141.358865490701, 0.3697181104890012, 0.3697181104890012),
  @("sCs", "Wnt6", "Fzd7", "sCs", 3, 1, 1.552337666666667, 4.657013, 0.9607189384054421, 0.960718938405442, 3, 1, 13.75539366666667, 41.266181, 0.5231823098264544, 0.5231823098264544, 21.35301570859478, 192.177141377353, 0.5026311532889783, 0.5026311532889783)
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $rowValues = $data[$i]
  $rowNum = $i + 2
  for ($j = 0; $j -lt $rowValues.Length; $j++) {
    $ws.Cells.Item($rowNum, $j + 1).Value = $rowValues[$j]
  }
}
